{"js": "// Apply the \"Analysis Project Update\" revision:\n//  - reword several bullet points under \"Insights from Initial Analysis\"\n//  - reword/merge the \"Upcoming Tasks\" bullets (3 bullets -> 1 bullet)\n//  - reword the closing paragraph\n//\n// Strategy: locate each paragraph by some stable (unchanged) anchor text via\n// Body.search(), then swap its wording by inserting a brand-new paragraph\n// (which inherits the anchor paragraph's list style / ilvl / numId) right\n// after it and deleting the original. This sidesteps leftover zero-width\n// markers (e.g. <w:proofErr/>) that a plain text replace can leave behind,\n// and works the same way whether the paragraph is a list item or plain text.\n\nasync function replaceParagraph(context, anchorText, newText) {\n  const body = context.document.body;\n  const results = body.search(anchorText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find paragraph starting with: \" + anchorText);\n  }\n  const para = results.items[0].paragraphs.getFirst();\n  para.insertParagraph(newText, Word.InsertLocation.after);\n  para.delete();\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 1) \"The most recent 30 day period has a low transaction count, ...\"\n// ---------------------------------------------------------------------\nawait replaceParagraph(\n  context,\n  \"The most recent 30 day period has a low transaction count\",\n  \"The most recent 30 day period has a relatively low transaction count (12), highlighting a potential issue with rewards transactions being reported in a timely manner.\"\n);\n\n// ---------------------------------------------------------------------\n// 2) \"...the top 5 brands from a transaction count standpoint were ...\"\n// ---------------------------------------------------------------------\nawait replaceParagraph(\n  context,\n  \"the top 5 brands from a transaction count standpoint were\",\n  \"In the month prior to the most recent, the top 5 brands from a transaction count standpoint were Amp, Sargento, Oscar Mayer, One, and Kraft.\"\n);\n\n// ---------------------------------------------------------------------\n// 3) \"Of those top 5 brands, the top 4 were not among the top 5 ...\"\n// ---------------------------------------------------------------------\nawait replaceParagraph(\n  context,\n  \"Of those top 5 brands\",\n  \"Despite Amp's surge overtaking Sargento, the top five brands have remained consistent month-over-month.\"\n);\n\n// ---------------------------------------------------------------------\n// 4) \"There was a 3X month-on-month increase ...\"\n// ---------------------------------------------------------------------\nawait replaceParagraph(\n  context,\n  \"There was a 3X\",\n  \"There has been a significant month-on-month (2X) increase in overall consumer transactions, potentially linked to the rise of Amp.\"\n);\n\n// ---------------------------------------------------------------------\n// 5) \"Digging deeper into the sharp increase in transactions ...\"\n// ---------------------------------------------------------------------\nawait replaceParagraph(\n  context,\n  \"deeper into the sharp increase in transactions\",\n  \"Delving deeper into the spike in transactions to identify underlying drivers and hypotheses, including exploring retailer and promotion data for additional context.\"\n);\n\n// ---------------------------------------------------------------------\n// 6) \"Investigating the last 30 days of transaction data ...\"\n// ---------------------------------------------------------------------\nawait replaceParagraph(\n  context,\n  \"Investigating the last 30 days of transaction data\",\n  \"Investigating the last 30 days of transaction data to uncover any potential issues in the feed, or to continue with our current analytical approach to accommodate the delay if it is an expected phenomenon (e.g. a delay in reporting to allow time for returns from customers).\"\n);\n\n// ---------------------------------------------------------------------\n// 7) Merge the 3 \"Collaborate with data partners ...\" bullets into 1.\n// ---------------------------------------------------------------------\n{\n  const body = context.document.body;\n\n  let results = body.search(\"Collaborate with data partners to address data ingestion\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find the 'Collaborate with data partners' paragraph\");\n  }\n  const para1 = results.items[0].paragraphs.getFirst();\n\n  results = body.search(\"Correcting invalid JSON files are repaired prior to transmittal.\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find the 'Correcting invalid JSON files' paragraph\");\n  }\n  const para2 = results.items[0].paragraphs.getFirst();\n\n  results = body.search(\"Create a more robust key between receipts and brands.\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find the 'Create a more robust key' paragraph\");\n  }\n  const para3 = results.items[0].paragraphs.getFirst();\n\n  // Insert the merged replacement paragraph right after the first bullet\n  // (it inherits para1's ListParagraph / ilvl=1 / numId=1 formatting), then\n  // remove all three original bullets.\n  para1.insertParagraph(\n    \"Collaborating with our data partners to address data ingestion issues, including rectifying invalid JSON files, enhancing key relationships between receipts and brands, scrubbing the list of brand codes, and considering the creation of an item-level table for improved granularity.\",\n    Word.InsertLocation.after\n  );\n  para1.delete();\n  para2.delete();\n  para3.delete();\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 8) Closing paragraph.\n// ---------------------------------------------------------------------\nawait replaceParagraph(\n  context,\n  \"Thank you for your ongoing support and involvement in this project.\",\n  \"Thank you for your ongoing support and involvement in this project. While we're encouraged by the insights uncovered thus far, we believe that further exploration will yield deeper insights to drive strategic decision-making in enhancing our customer service. Should you have any questions or wish to discuss the project in more detail, please don't hesitate to contact me directly.\"\n);\n", "ps1": "# Apply the \"Analysis Project Update\" revision:\n#  - reword several bullet points under \"Insights from Initial Analysis\"\n#  - reword/merge the \"Upcoming Tasks\" bullets (3 bullets -> 1 bullet)\n#  - reword the closing paragraph\n#\n# Strategy: locate each paragraph by some stable (unchanged) anchor text,\n# then swap its wording by inserting a brand-new paragraph right after it\n# (InsertParagraphAfter keeps the same list style / ilvl / numId as the\n# anchor paragraph) and deleting the original. This sidesteps leftover\n# zero-width markers (e.g. proofErr spell/grammar-check tags) that a plain\n# Range.Text assignment can leave behind when a run boundary falls in the\n# middle of such a marker.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphByText {\n    param($doc, [string]$needle)\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        if ($p.Range.Text.Contains($needle)) {\n            return $p\n        }\n    }\n    throw \"Paragraph containing '$needle' not found\"\n}\n\nfunction Set-ParagraphText {\n    param($doc, [string]$anchor, [string]$newText)\n    $p = Find-ParagraphByText $doc $anchor\n    $p.Range.InsertParagraphAfter()\n    $newIndex = $p.Index + 1\n    $newPara = $doc.Paragraphs.Item($newIndex)\n    $newPara.Range.Text = $newText\n    $p.Range.Delete()\n}\n\n# ---------------------------------------------------------------------\n# 1) \"The most recent 30 day period has a low transaction count, ...\"\n# ---------------------------------------------------------------------\nSet-ParagraphText $d \"The most recent 30 day period has a low transaction count\" `\n    \"The most recent 30 day period has a relatively low transaction count (12), highlighting a potential issue with rewards transactions being reported in a timely manner.\"\n\n# ---------------------------------------------------------------------\n# 2) \"...the top 5 brands from a transaction count standpoint were ...\"\n# ---------------------------------------------------------------------\nSet-ParagraphText $d \"the top 5 brands from a transaction count standpoint were\" `\n    \"In the month prior to the most recent, the top 5 brands from a transaction count standpoint were Amp, Sargento, Oscar Mayer, One, and Kraft.\"\n\n# ---------------------------------------------------------------------\n# 3) \"Of those top 5 brands, the top 4 were not among the top 5 ...\"\n# ---------------------------------------------------------------------\nSet-ParagraphText $d \"Of those top 5 brands\" `\n    \"Despite Amp's surge overtaking Sargento, the top five brands have remained consistent month-over-month.\"\n\n# ---------------------------------------------------------------------\n# 4) \"There was a 3X month-on-month increase ...\"\n# ---------------------------------------------------------------------\nSet-ParagraphText $d \"There was a 3X\" `\n    \"There has been a significant month-on-month (2X) increase in overall consumer transactions, potentially linked to the rise of Amp.\"\n\n# ---------------------------------------------------------------------\n# 5) \"Digging deeper into the sharp increase in transactions ...\"\n# ---------------------------------------------------------------------\nSet-ParagraphText $d \"deeper into the sharp increase in transactions\" `\n    \"Delving deeper into the spike in transactions to identify underlying drivers and hypotheses, including exploring retailer and promotion data for additional context.\"\n\n# ---------------------------------------------------------------------\n# 6) \"Investigating the last 30 days of transaction data ...\"\n# ---------------------------------------------------------------------\nSet-ParagraphText $d \"Investigating the last 30 days of transaction data\" `\n    \"Investigating the last 30 days of transaction data to uncover any potential issues in the feed, or to continue with our current analytical approach to accommodate the delay if it is an expected phenomenon (e.g. a delay in reporting to allow time for returns from customers).\"\n\n# ---------------------------------------------------------------------\n# 7) Merge the 3 \"Collaborate with data partners ...\" bullets into 1.\n# ---------------------------------------------------------------------\n# NOTE: Paragraph objects here are position-bound, not identity-bound, so\n# once InsertParagraphAfter() shifts everything below it down by one, any\n# previously-fetched Paragraph reference for those later bullets now\n# resolves to the wrong paragraph. To stay safe we (a) do the insertion\n# first, then (b) re-find each bullet to delete *immediately* before\n# deleting it, working from the bottom (highest index) up so that deleting\n# a later paragraph never invalidates the index of one we still need.\n$p1 = Find-ParagraphByText $d \"Collaborate with data partners to address data ingestion\"\n\n$p1.Range.InsertParagraphAfter()\n$newIndex = $p1.Index + 1\n$newPara = $d.Paragraphs.Item($newIndex)\n$newPara.Range.Text = \"Collaborating with our data partners to address data ingestion issues, including rectifying invalid JSON files, enhancing key relationships between receipts and brands, scrubbing the list of brand codes, and considering the creation of an item-level table for improved granularity.\"\n\n$p3 = Find-ParagraphByText $d \"Create a more robust key between receipts and brands.\"\n$p3.Range.Delete()\n\n$p2 = Find-ParagraphByText $d \"Correcting invalid JSON files are repaired prior to transmittal.\"\n$p2.Range.Delete()\n\n$p1again = Find-ParagraphByText $d \"Collaborate with data partners to address data ingestion\"\n$p1again.Range.Delete()\n\n# ---------------------------------------------------------------------\n# 8) Closing paragraph.\n# ---------------------------------------------------------------------\nSet-ParagraphText $d \"Thank you for your ongoing support and involvement in this project.\" `\n    \"Thank you for your ongoing support and involvement in this project. While we're encouraged by the insights uncovered thus far, we believe that further exploration will yield deeper insights to drive strategic decision-making in enhancing our customer service. Should you have any questions or wish to discuss the project in more detail, please don't hesitate to contact me directly.\"\n"}
